$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.58%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'4.479"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-5.69%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05872"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.75%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.95%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8518"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.10%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9279"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-2.15%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.05%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.04559"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'26.34%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07035"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.13%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03068"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.27%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-1.52%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001544"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.51%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006052"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-94.28%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006036"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.71%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.481"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.11%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.173"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.69%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D20").Value = "'0.3050"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.84%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1286"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.69%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.908"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.86%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04270"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.11%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.55%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004290"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-4.74%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-29.43%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'2.04%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03804"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.53%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-0.34%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.003887"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-37.42%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'5.62%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01385"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'26.26%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005379"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-2.24%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05391"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-50.60%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'10,961.08%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("E50").Style = "Normal"
